$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pak Elvis (Ustadz Elvis Indra) can no longer make it for his slot (row 15 / 21 Maret) -
# clear the assigned imam cell.
$ws.Range("D15").ClearContents()

# Row 26 (1 April / Ustadz Awan Karliawan) has now confirmed - mark the confirmation column.
$ws.Range("E26").Value2 = "✔️"

# Update the view so it reflects where the editor was looking (scrolled down, selection on E26).
$ws.Range("E26").Select()
$excel.ActiveWindow.ScrollRow = 10
